$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Molex Minifit")

# New row of data for the 8-pin Molex Micro-Fit header
$row = 30
$ws.Cells.Item($row, 1).Value = 430450815
$ws.Cells.Item($row, 2).Value = 8
$ws.Cells.Item($row, 3).Value = "3.0mm"
$ws.Cells.Item($row, 4).Value = "3.0mm"
$ws.Cells.Item($row, 5).Value = 90
$ws.Cells.Item($row, 6).Value = "No"
$ws.Cells.Item($row, 7).Value = "Gold"
$ws.Cells.Item($row, 8).Value = 430450815
$ws.Cells.Item($row, 9).Value = "Molex Micro-Fit.SchLib"
$ws.Cells.Item($row, 10).Value = 430450815
$ws.Cells.Item($row, 11).Value = "Molex Micro-Fit.PcbLib"

# Set the hyperlink cell (and its shared-string text) before the part-number
# cell so the new shared strings land in the same order as the source edit.
$linkCell = $ws.Cells.Item($row, 13)
$linkCell.Value = "https://www.digikey.com/en/products/detail/molex/0430450815/3044333"
$ws.Hyperlinks.Add($linkCell, "https://www.digikey.com/en/products/detail/molex/0430450815/3044333") | Out-Null

# Hyperlinks.Add stamps its own ad-hoc "hyperlink" format on the cell; re-apply
# the named Hyperlink cell style so M30 lands on the same shared style index
# (s="2") as the rest of column M instead of a brand-new one.
$linkCell.Style = "Hyperlink"

$ws.Cells.Item($row, 12).Value = "WM7067DKR-ND"

$ws.Range("J31").Select()
